$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the inline string values in column N for rows 2 through 304,
# leaving the cells present but empty (as inline strings with no content).
$ws.Range("N2:N304").Value = ""
